$wb = $excel.ActiveWorkbook

# --- Sheets ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1. Update status text for the 22e10bb2 row (row 7) everywhere it appears ---
# Before: "Ready for handoff"  -> After: "Handback transform failed"
$wsOverview.Range("E7").Value = "Handback transform failed"
$wsOverview.Range("F7").Value = "Handback transform failed"
$wsZhCn.Range("C7").Value = "Handback transform failed"
$wsDeDe.Range("C7").Value = "Handback transform failed"

# --- 2. Widen the "Error Detail" column (column P, the 16th column) on both
#        locale sheets so the new error text is readable ---
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17

# --- 3. Populate the new "Error Detail" messages for the 22e10bb2 row (row 7) ---
$wsZhCn.Range("P7").Value = "Handback file name: yy0stfww.nzj is different with handoff file name: 22e10bb2-2f82-4e2a-a45d-fa1ad852bec2.4db4096e6c40dda8feaa971f43d559186ce3a8b0.zh-cn."
$wsDeDe.Range("P7").Value = "Handback file name: yy0stfww.nzj is different with handoff file name: 22e10bb2-2f82-4e2a-a45d-fa1ad852bec2.4db4096e6c40dda8feaa971f43d559186ce3a8b0.de-de."
